$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 227, shifting rows 227-239 down to 228-240
$ws.Rows.Item(227).Insert()

# Populate the new row 227 with a fresh weekly data entry.
# Columns A-J, Q, R, T mirror the constant/contextual values used throughout
# this price-report block; D, K, L, M, N, O, P, S carry the new record values.
$ws.Cells.Item(227, 1).Value = 1
$ws.Cells.Item(227, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(227, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(227, 4).Value = 44706
$ws.Cells.Item(227, 5).Value = 15
$ws.Cells.Item(227, 6).Value = "Fruta"
$ws.Cells.Item(227, 7).Value = 100102
$ws.Cells.Item(227, 8).Value = "Cítricos"
$ws.Cells.Item(227, 9).Value = 100102003
$ws.Cells.Item(227, 10).Value = "Limón"
$ws.Cells.Item(227, 11).Value = "Sin especificar"
$ws.Cells.Item(227, 12).Value = "2a amarillo"
$ws.Cells.Item(227, 13).Value = 250
$ws.Cells.Item(227, 14).Value = 15000
$ws.Cells.Item(227, 15).Value = 16000
$ws.Cells.Item(227, 16).Value = 15500
$ws.Cells.Item(227, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(227, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(227, 19).Value = 775
$ws.Cells.Item(227, 20).Value = 20
